# Error Calculations and Plots
# Remove the "RM 232" and "SC 92" rows (rows that were entirely re-derived /
# dropped from the missing-data extract) and refresh the set of cells that
# were imputed/blanked for the remaining rows, per the authoritative edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two rows that disappear from the sheet -------------------
# Row 26 ("RM 232") is removed outright.
$ws.Rows.Item(26).Delete()
# After that shift, the row that used to hold "SC 92" is now row 27; remove it
# too, which brings the sheet down to A1:F33 (was A1:F35).
$ws.Rows.Item(27).Delete()

# --- 2. Blank out / fill in specific cells on the remaining rows -------------
# (Row numbers below refer to the CURRENT sheet layout, i.e. after the two
# row deletions above.)

# Row 5 ("RM 14"): column C (D-col) value removed
$ws.Cells.Item(5, 4).ClearContents()

# Row 11 ("RM 58"): column C (D-col) now has a value
$ws.Cells.Item(11, 4).Value = -15.5

# Row 19 ("RM 125"): column B (C-col) filled in, column C (D-col) blanked
$ws.Cells.Item(19, 3).Value = 13.2
$ws.Cells.Item(19, 4).ClearContents()

# Row 21 ("RM 135"): column B (C-col) blanked
$ws.Cells.Item(21, 3).ClearContents()

# Row 23 ("RM 140"): column B and C (C-col, D-col) filled in
$ws.Cells.Item(23, 3).Value = 12.2
$ws.Cells.Item(23, 4).Value = -13.9

# Row 25 ("RM 145"): column C (D-col) filled in
$ws.Cells.Item(25, 4).Value = -15.5

# Row 26 ("SC 5"): column A (B-col) blanked
$ws.Cells.Item(26, 2).ClearContents()

# Row 27 ("SC 101"): column A (B-col) filled in, columns B & C (C-col, D-col) blanked
$ws.Cells.Item(27, 2).Value = -20.4
$ws.Cells.Item(27, 3).ClearContents()
$ws.Cells.Item(27, 4).ClearContents()

# Row 29 ("SC 119"): column A (B-col) blanked, column C (D-col) blanked
$ws.Cells.Item(29, 2).ClearContents()
$ws.Cells.Item(29, 4).ClearContents()

# Row 33 ("SC 232"): columns B and C (C-col, D-col) filled in
$ws.Cells.Item(33, 3).Value = 10.4
$ws.Cells.Item(33, 4).Value = -14.1
